# RPAR_holdings.xlsx — "Add files via upload"
# Bumps the "as of" date in the confidential disclosure banner from
# 2021-04-27 to 2021-04-28, and refreshes the Weight/Percent Change
# figures (columns D/E, rows 2-15) to the newly re-run model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unlock it so the cells can be written,
# then restore protection once the edits are in place.
$ws.Unprotect()

$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."
# Re-fit the row height after the multi-line text write so it stays
# at the sheet's default (unchanged) height instead of an explicit one.
$ws.Rows("18:18").AutoFit()

# row -> (Weight, Percent Change)
$updates = @{
    2  = @(0.05776839152616171,  -0.0001376778338687545)
    3  = @(0.02371711934800563,   0.002738654147104835)
    4  = @(0.03188926560222285,   0.008043396932285729)
    5  = @(0.03086526187450464,   0.02855977152182776)
    6  = @(0.03721408498635759,   0.00854883522120109)
    7  = @(0.01923308987252559,   0.002289223762576675)
    8  = @(0.004879427471496127, -0.02037489812550941)
    9  = @(0.006945992923292504, -0.001526717557251867)
    10 = @(0.07018899337563704,   0.003399433427762277)
    11 = @(0.07030829478080809,   0.002262443438913797)
    12 = @(0.1470217494481332,    0.000649163300634692)
    13 = @(0.3851407263137151,    0.001316482359136506)
    14 = @(0.1148276024771399,    0.004363636363636347)
    15 = @(1,                     0.002948355276694903)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

$ws.Protect()
